# Deploying to gh-pages from @ Alvearie/alvearie-fhir-ig@8e4a450c507ef6f746e072652acbb72e9504f19a 🚀
#
# Updates the FHIR StructureDefinition workbook for "RelatedObservation":
#   - bump Version 5.0.0 -> 6.0.0
#   - refresh the generated Date
#   - fill in a real Publisher ("Alvearie Team")
#   - replace the stray duplicated "Contact" row with a "Jurisdiction" row
#     and drop the leftover duplicate row entirely
#   - give the root Extension row on the Elements sheet a real Short /
#     Definition instead of the generic placeholder text

$wb = $excel.ActiveWorkbook

# --- "Metadata" sheet ---
$metadata = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$metadata.Range("B3").Value = "6.0.0"

# Date: regenerated timestamp
$metadata.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher: was blank, now has a real value
$metadata.Range("B9").Value = "Alvearie Team"

# Row 10 used to be a duplicate "Contact" / "No display for ContactDetail"
# row; it becomes the new "Jurisdiction" row.
$metadata.Range("A10").Value = "Jurisdiction"
$metadata.Range("B10").Value = "United States of America"

# Row 11 was the other half of the duplicated "Contact" row - delete it so
# "Description" and everything below shifts up one row (A1:B21 -> A1:B20).
$metadata.Rows.Item(11).Delete()

# --- "Elements" sheet ---
$elements = $wb.Worksheets.Item("Elements")

# Root Extension row (row 2): Short / Definition get real content instead
# of the generic "Extension" / "An Extension" placeholders.
$elements.Range("K2").Value = "RelatedObservation"
$elements.Range("L2").Value = "Related observations that can be combined to fulfill a single time period for a given care gap or other measure population."
